$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "C2RobotStop.Ext"
$ws.Range("B8").Value = "C2RobotStopExt"
$ws.Range("C8").Value = "Extends robot stop"

$ws.Range("A9").Value = "C2RobotStop.Ret"
$ws.Range("B9").Value = "C2RobotStopRet"

$ws.Range("A16").Select()
